$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the last used row in column A (currently row 69) and append a new row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$dateCell = $ws.Cells.Item($newRow, 1)
$priceCell = $ws.Cells.Item($newRow, 2)

$dateCell.Value = "25-11-2025"
$priceCell.Value = "The price of gold in India today is ₹12,704 per gram for 24 karat gold, ₹11,645 per gram for 22 karat gold and ₹9,528 per gram for 18 karat gold (also called 999 gold)."

# Match formatting of the previous data row.
$ws.Cells.Item($lastRow, 1).Copy() | Out-Null
$dateCell.PasteSpecial(-4122) | Out-Null
$ws.Cells.Item($lastRow, 2).Copy() | Out-Null
$priceCell.PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
